$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark after "Takiego typu mogą być?"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Append the new paragraphs (with a fresh "_GoBack" bookmark wrapping the
#    final block) right before the end of the document body.
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:color w:val="5B9BD5" w:themeColor="accent1" /></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /><w:b /></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /><w:b /></w:rPr><w:t xml:space="preserve">Czym </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /><w:b /></w:rPr><w:t>jest nasza aplikacja:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>Nasza aplikacja została napisana z wykorzystaniem WPF i C#. Jest to uwarunkowane u podstaw tym, że nasze ostatnie lata na uczelni były związane właśnie z tym językiem i chcieliśmy nie tylko wykorzystać tę widzę i przetestować ale także poszerzyć tworząc coś własno co przysłuży się innym. Ale nie tylko to było powodem – postawiliśmy sobie za cel zrobić aplikację bazodanową opartą w pełni o usłu</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">gi lokalne dla instytucji muzealnych, które potrzebują aplikacji o niskim progu wejścia pod kątem szkoleń z obsługi, złożoności systemu i itp. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">System jest uniwersalny tzn. nie koniecznie musi znaleźć się w muzeum (muzeum jest przykładem z racji chęci wypełnienia luki systemem) ponieważ bez problemu może posłużyć np. bibliotece (książki zbiory), wypożyczalni (przedmiot i wypożyczający) i itp. Tutaj przechodzimy do drugiej dużej zalety </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>jest to system zaprojektowany w MVVM a to oznacza możliwość modyfikowania naszego projektu poprzez: bindowanie</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> w widoku i wykorzystanie już napisanych przez nas metod w modelu. Oczywiście nic nie stoi na przeszkodzie aby zmodyfikować nieco projekt pod klienta preferencje jednak fundament pozostaje taki sam.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">Warto zwrócić uwagę na zasadność użycia </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>SQLite</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> do naszej pracy – bo oprócz oczywistych zalet jak lekkość, szybkość</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> ważnym elementem jest też dla nas możliwość pliku .</w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>db</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>. Co to oznacza? Otóż w małych aplikacjach</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> baza jest wystarczająca</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> jednak gdyby zaszła potrzeba </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>przeniesienia bazy do innego serwisu czy po prostu archiwizacja w chmurze to jest to bardzo proste do przeprowadzenia z racji ustandaryzowania tego formatu u jego podstaw. Zatem nie ograniczamy administratora tego systemu a wręcz ułatwiamy mu modyfikację danych i ich moderowanie. Oczywiście hasła w bazie są przed nim ukryte funkcją hasz</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> SHA256 a sama baza gdyby utworzyć kopię zabezpieczona kluczem symetrycznym AES.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">Wyjaśnienie tabel specjalnych: </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>current_user</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> to tabela celowo przechowująca sam login obecnego </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>usera</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>. Jest tak ze względów bezpieczeństwa aby dane użytkownika ale też ich obróbka była prosta w opanowaniu – bowiem do sytuacji takich jak wysyłanie maili czy rejestr zdarzeń w tabeli historia – jest potrzebny tylko login. Czyści się po zamknięciu.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:lastRenderedPageBreak /><w:t xml:space="preserve">Tabela historia jest oparta wyłącznie na </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>triggerach</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> czyli wyzwalaczach, które kiedy coś zadzieje się w bazie to wykonają (w tym przypadku) rejestr o zmianie w bazie – służy głównie administratorowi jako element jego uprawnień.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto" /><w:ind w:left="708" /><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">Optymalizacja aplikacji – polegała głównie na zwolnieniu zasobów (połączeniu </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>DbConnect</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> i </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>Window</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">) za pomocą funkcji </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>Dispose</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">() i Close() z wykorzystaniem </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve">instrukcji </w:t></w:r><w:proofErr w:type="spellStart" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>using</w:t></w:r><w:proofErr w:type="spellEnd" /><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t>. Pozwoliło to nam na większą kontrolę zasobów systemu. Natomiast samo wykorzystanie WPF, który wykorzystuje grafikę wektorową zamiast rastrowej pozwoliło nam na zmniejszenie wymagań aplikacji w tym przede wszystkim wykorzystanie ram-u na poziomie ~180 MB. Zrzucenie</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" /></w:rPr><w:t xml:space="preserve"> rejestrowania na bazę a przy tym przyspieszanie procesu rejestracji zdarzeń.</w:t></w:r><w:bookmarkEnd w:id="0" /></w:p>
'@

$insertionPoint.InsertXML($newParagraphsXml)
